$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D27").Value = "RLHF 외에 LLM이 피드백을 학습할 수 있는 방법은 무엇이 있을까?"
$ws.Range("E27").Value = "https://tech.scatterlab.co.kr/alt-rlhf/"

$ws.Range("D36").Value = "Multimodal Contrastive learning with various data domains"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/430"

$ws.Range("D51").Value = "[aws] EC2 인스턴스 AMI 생성하는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/aws-EC2-%EC%9D%B8%EC%8A%A4%ED%84%B4%EC%8A%A4-AMI-%EC%83%9D%EC%84%B1%ED%95%98%EB%8A%94-%EB%B0%A9%EB%B2%95"

$wb.Save()
